$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the cell value (creates/uses a shared string entry, matching the
# target workbook's xl/sharedStrings.xml + xl/worksheets/sheet1.xml).
$ws.Range("K8").Value = "asdasdasdvxzv"

# Move the selection onto K8, matching the <selection activeCell="K8"
# sqref="K8"/> recorded in the sheet view.
[void]$ws.Range("K8").Select()
